$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2021-09-07"

# Update header cell text for September column
$ws.Range("B1").Value = "September 2021 (through September 07)"

# Update / add data values for 2021-09-15 (per commit message) across
# several neighborhoods in the September 2021 column (B) and various
# historical month columns.
$ws.Range("B2").Value = 5      # Garfield Park
$ws.Range("K4").Value = 1      # Humboldt Park
$ws.Range("K6").Value = 2      # Roseland
$ws.Range("K9").Value = 1      # Little Village
$ws.Range("T9").Value = 1      # Little Village
$ws.Range("AL10").Value = 1    # West Town
$ws.Range("AU10").Value = 3    # West Town
$ws.Range("T13").Value = 1     # Chatham
$ws.Range("BD16").Value = 1    # West Pullman
$ws.Range("B17").Value = 1     # South Shore
$ws.Range("AL18").Value = 1    # Grand Boulevard
$ws.Range("BD20").Value = 2    # Englewood
$ws.Range("T22").Value = 1     # South Chicago
$ws.Range("B32").Value = 1     # Chicago Lawn
$ws.Range("T91").Value = 3     # Pullman
$ws.Range("AU98").Value = 1    # West Lawn
$ws.Range("AL99").Value = 2    # West Ridge
